$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F14").Value = 9
$ws.Range("G14").Value = 12458.79
$ws.Range("B15").Value = 78339.50999999999
$ws.Range("F32").Value = 11
$ws.Range("G32").Value = 281.71
$ws.Range("F36").Value = 110
$ws.Range("G36").Value = 3099.8
$ws.Range("B40").Value = 52784.43
$ws.Range("F42").Value = 67
$ws.Range("G42").Value = 13183.59
$ws.Range("F54").Value = 216
$ws.Range("G54").Value = 12117.6
$ws.Range("F60").Value = 46
$ws.Range("G60").Value = 2051.6
$ws.Range("F61").Value = 105
$ws.Range("G61").Value = 5854.8
$ws.Range("F62").Value = 37
$ws.Range("G62").Value = 825.84
$ws.Range("F67").Value = 192
$ws.Range("G67").Value = 50060.16
$ws.Range("B72").Value = 173644.06
$ws.Range("F101").Value = 1
$ws.Range("G101").Value = 371.8
$ws.Range("B103").Value = 14911.13
$ws.Range("F118").Value = 33
$ws.Range("G118").Value = 2606.34
$ws.Range("F119").Value = 103
$ws.Range("G119").Value = 14456.05
$ws.Range("B129").Value = 67277.09
$ws.Range("B132").Value = 64196
$ws.Range("B133").Value = 65258
$ws.Range("F145").Value = 26
$ws.Range("G145").Value = 1771.38
$ws.Range("B153").Value = 19118.52
$ws.Range("F158").Value = 127
$ws.Range("G158").Value = 6059.17
$ws.Range("F160").Value = 287
$ws.Range("G160").Value = 9568.58
$ws.Range("B161").Value = 33646.06
$ws.Range("F190").Value = 45
$ws.Range("G190").Value = 3690
$ws.Range("B199").Value = 55514.97
$ws.Range("F212").Value = 49
$ws.Range("G212").Value = 3175.2
$ws.Range("B214").Value = 3175.2
$ws.Range("F228").Value = 288
$ws.Range("G228").Value = 5328
$ws.Range("F233").Value = 22
$ws.Range("G233").Value = 2521.2
$ws.Range("B235").Value = 12713.67
$ws.Range("F284").Value = 26
$ws.Range("G284").Value = 3524.04
$ws.Range("B296").Value = 64983
$ws.Range("C296").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F296").Value = 6
$ws.Range("G296").Value = 514.08
$ws.Range("B297").Value = 66194
$ws.Range("C297").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F297").Value = 17
$ws.Range("G297").Value = 1456.56
$ws.Range("B298").Value = 64985
$ws.Range("C298").Value = "HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S"
$ws.Range("F298").Value = 12
$ws.Range("G298").Value = 1052.4
$ws.Range("B299").Value = 66196
$ws.Range("C299").Value = "HIM-Total Care Baby Pants Drapers-Xl-9S"
$ws.Range("F299").Value = 1
$ws.Range("G299").Value = 87.7
$ws.Range("B301").Value = 95388.81
$ws.Range("B310").Value = 55373
$ws.Range("E310").Value = 163.62
$ws.Range("F310").Value = -94
$ws.Range("G310").Value = -13562.32
$ws.Range("B311").Value = 63520
$ws.Range("E311").Value = 153.4
$ws.Range("F311").Value = 35
$ws.Range("G311").Value = 5049.8
$ws.Range("B312").Value = 57802
$ws.Range("E312").Value = 162.71
$ws.Range("F312").Value = -79
$ws.Range("G312").Value = -11334.92
$ws.Range("B313").Value = 63531
$ws.Range("E313").Value = 152.53
$ws.Range("F313").Value = 24
$ws.Range("G313").Value = 3443.52
$ws.Range("B314").Value = 55356
$ws.Range("E314").Value = 54.04
$ws.Range("F314").Value = -158
$ws.Range("G314").Value = -7527.12
$ws.Range("B315").Value = 63510
$ws.Range("E315").Value = 50.66
$ws.Range("F315").Value = 74
$ws.Range("G315").Value = 3525.36
$ws.Range("B323").Value = 63560
$ws.Range("E323").Value = 134.87
$ws.Range("F323").Value = 1
$ws.Range("G323").Value = 126.86
$ws.Range("B324").Value = 60325
$ws.Range("E324").Value = 151.57
$ws.Range("F324").Value = -102
$ws.Range("G324").Value = -12939.72
$ws.Range("F367").Value = 192
$ws.Range("G367").Value = 26993.28
$ws.Range("B369").Value = 58919.66
$ws.Range("F371").Value = 4
$ws.Range("G371").Value = 221.32
$ws.Range("F372").Value = 42
$ws.Range("G372").Value = 2323.86
$ws.Range("F376").Value = 160
$ws.Range("G376").Value = 26558.4
$ws.Range("F377").Value = 61
$ws.Range("G377").Value = 9166.469999999999
$ws.Range("B378").Value = 47546.15
$ws.Range("F403").Value = 56
$ws.Range("G403").Value = 2025.52
$ws.Range("F408").Value = 15
$ws.Range("G408").Value = 514.65
$ws.Range("F409").Value = 60
$ws.Range("G409").Value = 2432.4
$ws.Range("F414").Value = 167
$ws.Range("G414").Value = 2646.95
$ws.Range("F417").Value = 112
$ws.Range("G417").Value = 3547.04
$ws.Range("B423").Value = 155801.71
$ws.Range("F438").Value = 52
$ws.Range("G438").Value = 2517.32
$ws.Range("B444").Value = 20390.34
$ws.Range("F517").Value = 172
$ws.Range("G517").Value = 17177.64
$ws.Range("F519").Value = 208
$ws.Range("G519").Value = 10289.76
$ws.Range("F523").Value = 132
$ws.Range("G523").Value = 7823.64
$ws.Range("B531").Value = 107267.61
$ws.Range("F534").Value = 227
$ws.Range("G534").Value = 3600.22
$ws.Range("B541").Value = 18821.22
$ws.Range("F550").Value = 26
$ws.Range("G550").Value = 1609.4
$ws.Range("B562").Value = 34892.76
$ws.Range("F607").Value = 8
$ws.Range("G607").Value = 1036.24
$ws.Range("B609").Value = 6131.29
$ws.Range("F611").Value = 155
$ws.Range("G611").Value = 20630.5
$ws.Range("B613").Value = 20630.5
$ws.Range("F618").Value = 217
$ws.Range("G618").Value = 32638.97
$ws.Range("F623").Value = 2
$ws.Range("G623").Value = 96.23999999999999
$ws.Range("F626").Value = 345
$ws.Range("G626").Value = 27113.55
$ws.Range("F631").Value = 288
$ws.Range("G631").Value = 10607.04
$ws.Range("B634").Value = 191071.69
$ws.Range("F680").Value = 372
$ws.Range("G680").Value = 60676.92
$ws.Range("B686").Value = 61689.47
$ws.Range("F721").Value = 122
$ws.Range("G721").Value = 12629.44
$ws.Range("B723").Value = 19484.39
$ws.Range("B724").Value = 2199797.37
$ws.Range("B725").Value = 2199797.37

Write-Output "Applied 163 cell changes"